$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 171860
$ws.Range("C4").Value = 162661
$ws.Range("C7").Value = 5.35
$ws.Range("C8").Value = 65.95999999999999
